$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2;  B="5478"; C="3725"},
    @{Row=3;  B="9803"; C="5439"},
    @{Row=4;  B="9494"; C="5942"},
    @{Row=5;  B="3662"; C="7821"},
    @{Row=6;  B="7664"; C="6117"},
    @{Row=7;  B="1025"; C="4573"},
    @{Row=8;  B="9554"; C="6223"},
    @{Row=9;  B="4942"; C="6392"},
    @{Row=10; B="4129"; C="231"},
    @{Row=11; B="5682"; C="104"},
    @{Row=12; B="9966"; C="9038"},
    @{Row=13; B="1065"; C="7508"},
    @{Row=14; B="1577"; C="851"},
    @{Row=15; B="8852"; C="2227"},
    @{Row=16; B="4488"; C="2903"},
    @{Row=17; B="1920"; C="4206"},
    @{Row=18; B="5124"; C="8115"},
    @{Row=19; B="4176"; C="5490"},
    @{Row=20; B="5962"; C="6222"},
    @{Row=21; B="8063"; C="269"},
    @{Row=22; B="1288"; C="954"},
    @{Row=23; B="8391"; C="237"},
    @{Row=24; B="6891"; C="3956"},
    @{Row=25; B="8532"; C="8574"},
    @{Row=26; B="8677"; C="6796"}
)

foreach ($item in $data) {
    $r = $item.Row
    $b = $item.B
    $c = $item.C
    $ws.Cells.Item($r, 2).NumberFormat = "@"
    $ws.Cells.Item($r, 3).NumberFormat = "@"
    $ws.Cells.Item($r, 4).NumberFormat = "@"
    $ws.Cells.Item($r, 2).Value = $b
    $ws.Cells.Item($r, 3).Value = $c
    $ws.Cells.Item($r, 4).Value = "$b-$c"
}
